{"js": "const body = context.document.body;\n\nconst pairs = [\n  [\"393\u00d73=\", \"297\u00d77=\"],\n  [\"693\u00d72=\", \"402\u00d78=\"],\n  [\"314\u00d73=\", \"503\u00d75=\"],\n  [\"607\u00d78=\", \"294\u00d77=\"],\n  [\"505\u00d78=\", \"883\u00d78=\"],\n  [\"202\u00d72=\", \"160\u00d72=\"],\n  [\"503\u00d77=\", \"511\u00d76=\"],\n  [\"816\u00d74=\", \"656\u00d78=\"],\n  [\"722\u00d73=\", \"335\u00d77=\"],\n  [\"947\u00d76=\", \"498\u00d72=\"],\n  [\"925\u00d76=\", \"878\u00d76=\"],\n  [\"160\u00d74=\", \"272\u00d75=\"],\n  [\"342\u00d72=\", \"145\u00d78=\"],\n  [\"681\u00d75=\", \"743\u00d78=\"],\n  [\"780\u00d74=\", \"815\u00d74=\"],\n  [\"962\u00d72=\", \"954\u00d72=\"],\n  [\"944\u00d72=\", \"155\u00d79=\"],\n  [\"307\u00d76=\", \"729\u00d76=\"],\n  [\"583\u00d76=\", \"605\u00d79=\"],\n  [\"109\u00d79=\", \"467\u00d77=\"],\n  [\"743\u00d77=\", \"497\u00d78=\"],\n  [\"118\u00d77=\", \"722\u00d78=\"],\n  [\"797\u00d75=\", \"343\u00d72=\"],\n  [\"322\u00d72=\", \"227\u00d74=\"],\n  [\"306\u00d79=\", \"455\u00d75=\"],\n];\n\nconst allResults = [];\nfor (const [find, replace] of pairs) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items/text');\n  allResults.push({ results, replace });\n}\nawait context.sync();\n\nfor (const { results, replace } of allResults) {\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('393\u00d73=', '297\u00d77='),\n    @('693\u00d72=', '402\u00d78='),\n    @('314\u00d73=', '503\u00d75='),\n    @('607\u00d78=', '294\u00d77='),\n    @('505\u00d78=', '883\u00d78='),\n    @('202\u00d72=', '160\u00d72='),\n    @('503\u00d77=', '511\u00d76='),\n    @('816\u00d74=', '656\u00d78='),\n    @('722\u00d73=', '335\u00d77='),\n    @('947\u00d76=', '498\u00d72='),\n    @('925\u00d76=', '878\u00d76='),\n    @('160\u00d74=', '272\u00d75='),\n    @('342\u00d72=', '145\u00d78='),\n    @('681\u00d75=', '743\u00d78='),\n    @('780\u00d74=', '815\u00d74='),\n    @('962\u00d72=', '954\u00d72='),\n    @('944\u00d72=', '155\u00d79='),\n    @('307\u00d76=', '729\u00d76='),\n    @('583\u00d76=', '605\u00d79='),\n    @('109\u00d79=', '467\u00d77='),\n    @('743\u00d77=', '497\u00d78='),\n    @('118\u00d77=', '722\u00d78='),\n    @('797\u00d75=', '343\u00d72='),\n    @('322\u00d72=', '227\u00d74='),\n    @('306\u00d79=', '455\u00d75=')\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
